$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "fecha"/"hora generado" labels in J2:K3 ---
$ws.Range("J2").Value = "Fecha en que se generó el reporte: "

# K2 must stay plain text "2023/05/16" (not auto-converted to a date serial)
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "2023/05/16"
$ws.Range("K2").NumberFormat = "general"

$ws.Range("J3").Value = "Hora en que se generó el reporte: "
$ws.Range("K3").Value = " 15:43"

# --- Update price/total figures ---
$ws.Range("F9").Value = 127.35
$ws.Range("F10").Value = 2028.25

# --- Column widths for the new J/K columns ---
$ws.Columns("J").ColumnWidth = 32.140625
$ws.Columns("K").ColumnWidth = 10.7109375
